$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at "L". This pushes the existing "Feed" column
# (and everything to its right, including the wide footer/sidebar bar
# on row 10) one column to the right, while the new blank column
# inherits the formatting of the old column L (so it lines up with the
# existing K:L column-width/style definition and the J1:L1 merged
# "Week 1" header grows to J1:M1 automatically).
$ws.Columns("L").Insert()

# Fill in the new "Egg Weight" column.
$ws.Range("L2").Value = "Egg Weight"
$ws.Range("L3").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("L5").Value = 0

# Restore the cursor/selection position recorded in the sheet.
[void]$ws.Range("M13").Select()
